$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new column before the current column B ("agent_num"), shifting
# agent_num and everything after it one column to the right.
$ws1.Columns.Item(2).Insert()

# Populate the new "number_of_run" column (header + data rows).
$ws1.Range("B1").Value = "number_of_run"
$ws1.Range("B2").Value = 1
$ws1.Range("B3").Value = 2
$ws1.Range("B4").Value = 1

# Match the new column's width to the authored width (~14 chars, best-fit).
$ws1.Columns.Item(2).ColumnWidth = 13.29

# Give the "scenarios" sheet a print/page setup definition.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Make "scenarios" the active sheet/tab, with A5 selected.
$ws1.Activate()
$ws1.Range("A5").Select() | Out-Null
